$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 689462
$ws.Range("C2").Value = 1.84
$ws.Range("D2").Value = 0.73
$ws.Range("I2").Value = 28.38

# Row 3
$ws.Range("B3").Value = 689462
$ws.Range("C3").Value = 2.32
$ws.Range("D3").Value = 0.29
$ws.Range("I3").Value = 51.36

# Row 4
$ws.Range("B4").Value = 689462
$ws.Range("C4").Value = 2.49
$ws.Range("D4").Value = 0.54
$ws.Range("I4").Value = 123.82

# Row 5
$ws.Range("B5").Value = 689462
$ws.Range("C5").Value = 2.72
$ws.Range("D5").Value = 0.97
$ws.Range("I5").Value = 205.83

# Row 6
$ws.Range("B6").Value = 689462
$ws.Range("C6").Value = 2.48
$ws.Range("D6").Value = 0.52
$ws.Range("I6").Value = 92.42

# Row 7
$ws.Range("B7").Value = 689462
$ws.Range("C7").Value = 1.97
$ws.Range("D7").Value = 0.58
$ws.Range("I7").Value = 27.22
